$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 field updates per diff
$ws.Range("A2").Value = 102077473
$ws.Range("B2").Value = 96367
$ws.Range("E2").Value = 219874
$ws.Range("F2").Value = "Nattviol"
$ws.Range("G2").Value = "Platanthera bifolia"
$ws.Range("H2").Value = "(L.) Rich."
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "3"

# Newly-present (empty) cells in the diff
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("AF2").Value = ""

$ws.Range("P2").Value = "Tvetaspåret, Tveta, Srm"
$ws.Range("Q2").Value = 647720.9098417715
$ws.Range("R2").Value = 6560694.968483768
$ws.Range("S2").Value = 10

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2022-06-28"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2022-07-05"

# AI2 ("barrskog") removed entirely
$ws.Range("AI2").ClearContents()

$ws.Range("AW2").Value = "Åsa Johansson"
$ws.Range("AX2").Value = "Åsa Johansson"
